# "adding import Excel ListeAdmis"
#
# This recreates (via the Excel object model) a workbook edit that:
#   1. Duplicates the existing "Feuil1" sheet into a new sheet named "gege"
#      (this is how the new sheet's identical header/ID/Nom/prenom rows +
#      styling/formatting come to exist).
#   2. On the original "Feuil1" sheet, the last row's "prenom" entry stops
#      being the text "Autre" and becomes the number 9.
#   3. On the new "gege" sheet, the last row's "Nom"/"prenom" entries become
#      "VAOVAO" / "DJ" respectively.
#   4. The new "gege" tab ends up as the active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) "Feuil1": row 4, column C ("Autre") turns into the number 9 --------
$ws1.Range("C4").Value = 9

# --- 2) Duplicate "Feuil1" right after itself to create the new sheet -----
$ws1.Copy($null, $ws1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "gege"

# --- 3) Update the new sheet's last row (order matters for shared-string
#        ordering: C4 "DJ" must be registered before B4 "VAOVAO") ----------
$ws2.Range("C4").Value = "DJ"
$ws2.Range("B4").Value = "VAOVAO"

# --- 4) Leave "Feuil1" showing the A1:C5 selection, then activate "gege" --
$ws1.Select()
$ws1.Range("A1:C5").Select()
$ws2.Select()
